$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new column C with English translations of column B,
# in the same order the rows appear (matches shared-string append order).
$ws.Range("C3").Value   = "First case of the virus in Chile"
$ws.Range("C11").Value  = "Pandemic declared"
$ws.Range("C14").Value  = "Phase 3 "
$ws.Range("C16").Value  = "Phase 4. Influenza vaccination campaign started in advance. "
$ws.Range("C18").Value  = "Closing of national borders"
$ws.Range("C19").Value  = "State of constitutional emergency. Closure of the malls throughout Chile. Use of 2% constitutional (4.7% of GDP)"
$ws.Range("C20").Value  = "Sanitary customs in eight regions of Chile "
$ws.Range("C21").Value  = "First coronavirus death in Chile"
$ws.Range("C22").Value  = "Curfew throughout the country."
$ws.Range("C25").Value  = "Total quarantine in seven communes of the RM"
$ws.Range("C31").Value  = "Mandatory quarantine extended for one week in six of the seven communes of the MR "
$ws.Range("C36").Value  = "Call to wear masks in public places"
$ws.Range("C39").Value  = "Mask use on public and private transport is mandatory."
$ws.Range("C40").Value  = "The western sector of Puente Alto commune (west of Concha y Toro Avenue) is added to the quarantine."
$ws.Range("C44").Value  = "Quarantine ends in some communes of Stgo; Chillán and Chillán Viejo; and Hualpén and San Pedro de la Paz."
$ws.Range("C50").Value  = "It is surpassed the ten thousand total infections at national level.  President Sebastián Piñera calls for the reactivation of the economy and instructs the gradual return of public officials to their jobs in a face-to-face manner."
$ws.Range("C54").Value  = "Suspension of the return to classes indefinitely"
$ws.Range("C69").Value  = " 12 communes are quarantined. Ñuñoa ends its confinement."
$ws.Range("C74").Value  = "Minsal reduces from 14 to 4 days the maximum license for suspected cases of COVID-19"
$ws.Range("C76").Value  = "R.M. goes into quarantine"
$ws.Range("C89").Value  = "Chile has surpassed mainland China in the total number of officially registered infections."
$ws.Range("C100").Value = "New counting methodology."
$ws.Range("C102").Value = "Quarantine for Valparaiso, Viña del Mar and 6 other communes "
$ws.Range("C105").Value = "Jaime Mañalich leaves office and Enrique Paris takes over. "

# Header for the new column, added last.
$ws.Range("C1").Value = "Hito_eng"

# Leave the selection on the new header cell, as in the edited workbook.
$ws.Range("C1").Select()
